$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 341, shifting all existing rows
# (341-412) down to (342-413).
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new data point.
$ws.Cells.Item(341, 1).Value = 10
$ws.Cells.Item(341, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(341, 3).Value = "La Araucanía"
$ws.Cells.Item(341, 4).Value = 44995
$ws.Cells.Item(341, 5).Value = 9
$ws.Cells.Item(341, 6).Value = 100112001
$ws.Cells.Item(341, 7).Value = "Berenjena"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 35
$ws.Cells.Item(341, 11).Value = 14000
$ws.Cells.Item(341, 12).Value = 14000
$ws.Cells.Item(341, 13).Value = 14000
$ws.Cells.Item(341, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(341, 15).Value = "Región del Maule"
$ws.Cells.Item(341, 16).Value = 350
$ws.Cells.Item(341, 17).Value = 40
$ws.Cells.Item(341, 18).Value = "Hortaliza"
